$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers: F1 = in_degree, G1 = name, H1 = out_degree
$ws.Range("F1").Value = "in_degree"
$ws.Range("G1").Value = "name"
$ws.Range("H1").Value = "out_degree"

# Match the existing bold/border/centered header formatting used by B1:E1
$ws.Range("B1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row data: row, in_degree, name, out_degree
$data = @(
    @(2, 26, "Ahsoka Tano", 25),
    @(3, 14, "Rex", 18),
    @(4, 0, "Ridge", 0),
    @(5, 3, "4-A7", 2),
    @(6, 1, "Wurtz", 1),
    @(7, 16, "R2-D2", 20),
    @(8, 10, "Luminara Unduli", 8),
    @(9, 1, "Unidentified Advanced Recon Force trooper lieutenant", 5),
    @(10, 27, "Darth Sidious", 20),
    @(11, 1, "Matchstick", 1),
    @(12, 13, "Qui-Gon Jinn", 9),
    @(13, 8, "Whorm Loathsom", 2),
    @(14, 3, "Kharrus", 2),
    @(15, 22, "Padmé Amidala", 19),
    @(16, 21, "Grievous", 14),
    @(17, 6, "Onaconda Farr", 6),
    @(18, 0, "Unidentified Clone Sergeant", 0),
    @(19, 13, "C-3PO", 18),
    @(20, 24, "Jabba Desilijic Tiure", 13),
    @(21, 14, "Ziro Desilijic Tiure/Canon", 10),
    @(22, 7, "Wat Tambor", 9),
    @(23, 20, "Mace Windu", 20),
    @(24, 1, "Unidentified Advanced Recon Force trooper commander", 5),
    @(25, 11, "Wullf Yularen", 12),
    @(26, 0, "Unidentified clone trooper pilot lieutenant", 2),
    @(27, 0, "R4-P17", 5),
    @(28, 0, "Gold Six (clone trooper)", 0),
    @(29, 1, "KRONOS-327", 1),
    @(30, 21, "Asajj Ventress", 19),
    @(31, 12, "Plo Koon", 16),
    @(32, 1, "Broadside", 5),
    @(33, 2, "CC-1010", 10),
    @(34, 1, "CC-2237", 6),
    @(35, 0, "Unidentified Nikto guard (Jabba's Palace)", 3),
    @(36, 0, "R2-C2", 1),
    @(37, 10, "Kit Fisto", 15),
    @(38, 30, "Anakin Skywalker", 28),
    @(39, 13, "Nute Gunray", 15),
    @(40, 8, "Shmi Skywalker Lars", 7),
    @(41, 0, "R2-KT", 10),
    @(42, 24, "Dooku", 15),
    @(43, 3, "TC-70", 2),
    @(44, 8, "CC-2224", 15),
    @(45, 1, "TB-2", 1),
    @(46, 24, "Obi-Wan Kenobi", 25),
    @(47, 17, "Yoda", 17),
    @(48, 23, "Rotta", 4)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 6).Value = $item[1]
    $ws.Cells.Item($r, 7).Value = $item[2]
    $ws.Cells.Item($r, 8).Value = $item[3]
}
